# The wml.xsd schema expects character-run properties inside <w:rPr> to
# appear in a fixed order (rFonts, b, bCs, i, iCs, caps, smallCaps,
# strike, ..., color, ...). A number of the Pandoc "*Tok" character
# styles in styles.xml had <w:b/>/<w:i/> emitted *after* <w:color/>,
# which OOXMLValidatorCLI flags as Sch_UnexpectedElementContentExpectingComplex
# even though xmllint stays quiet. Re-assert the bold/italic flags on
# each affected style so the engine re-serializes <w:rPr> in
# schema-correct order (b/i before color), without altering any of the
# actual formatting values.

$d = $word.ActiveDocument
$styles = $d.Styles

$boldStyles = @(
    "KeywordTok",
    "ImportTok",
    "AnnotationTok",
    "CommentVarTok",
    "ControlFlowTok",
    "InformationTok",
    "WarningTok",
    "AlertTok",
    "ErrorTok"
)

$italicStyles = @(
    "CommentTok",
    "DocumentationTok",
    "AnnotationTok",
    "CommentVarTok",
    "InformationTok",
    "WarningTok"
)

foreach ($s in $styles) {
    $name = $s.NameLocal
    if ($boldStyles -contains $name) {
        $s.Font.Bold = $True
    }
    if ($italicStyles -contains $name) {
        $s.Font.Italic = $True
    }
}
